# Insert a new column before BO ("nom"), shifting "nom" -> BP and
# "url_produit" -> BQ, to make room for a new price-history snapshot
# column, then populate it:
#  - header (row 1) gets the new timestamp
#  - data rows that already have a tracked price in BN (rows 2-80) get
#    that same price carried forward into the new BO column
#  - data rows with no tracked price yet (rows 81-206) are left blank

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("BO:BO").Insert()

$ws.Range("BO1").Value = "2026-01-30 17:24:55"

$lastRow = 206
for ($r = 2; $r -le $lastRow; $r++) {
    $priceCell = $ws.Cells.Item($r, 66)   # column BN
    $price = $priceCell.Value2()
    if ($price -ne $null -and $price -ne "") {
        $ws.Cells.Item($r, 67).Value = $price   # column BO
    }
}
